$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.493.77'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.847.80'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9993'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6302'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -1.82%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2907'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.91'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07741'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.849.69'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001024'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.316'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.526.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '229.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.518'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.13%  '
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.509'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1361'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.46%  '
$ws.Range('E27').Value = '  -0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06602'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +16.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.463'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.489'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.096'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.848'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.41%  '
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6963'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.64%  '
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01870'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.832'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.252.17'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.783'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9388'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.03%  '
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.008.82'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.25'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.085'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.721'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1156'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.024'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3929'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.00000000111'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.32%  '
